# Implements the (unfinished) "Penalty Reward System" edit:
# - "Weekly Quantity" sheet: rows 3 & 4 get replaced with what used to be
#   rows 7 & 8 (i.e. weeks 45158.99999999999/20 and 45179.99999999999/20),
#   and the now-redundant trailing rows 5-8 are removed.
# - "Monthly Trend" sheet: row 3 gets replaced with what used to be row 4
#   (45169.99999999999/20), row 4 becomes the old row 5
#   (45199.99999999999/20), and the now-redundant trailing row 5 is removed.

$wb = $excel.ActiveWorkbook

# --- Weekly Quantity ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$ws1.Cells.Item(3, 1).Value = 45158.99999999999
$ws1.Cells.Item(3, 2).Value = 20

$ws1.Cells.Item(4, 1).Value = 45179.99999999999
$ws1.Cells.Item(4, 2).Value = 20

# Remove the old rows 5-8, which are no longer needed.
$ws1.Rows.Item(5).Resize(4).Delete()

# --- Monthly Trend ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(3, 1).Value = 45169.99999999999
$ws2.Cells.Item(3, 2).Value = 20

$ws2.Cells.Item(4, 1).Value = 45199.99999999999
$ws2.Cells.Item(4, 2).Value = 20

# Remove the old row 5, which is no longer needed.
$ws2.Rows.Item(5).Delete()
